# Generate Report for Handoff
#
# A new handoff cycle ran for e2e\b.md: a fresh target xliff
# (b.*.xlf) was generated for both zh-cn and de-de, the handback for
# that file is not yet in sync with the newest source revision (so its
# status flips from "Handed back: in sync with en-US" to "Ready for
# handoff" and Content Duplicate flips from True to False), and an
# Error Detail message explains the stale handback. The Overview sheet
# mirrors the new status/date for b.md. The Error Detail column is
# widened so the long message is readable.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4942b1b7337f8e45c9868b6e98a887feb5e6dcb7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93dba24b6072debb9715407ef90b496a253427ab/e2e/b.md."

# Helper: write a literal "True"/"False"-looking string as TEXT instead
# of letting Excel auto-coerce it to a Boolean cell.
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $range.Worksheet.Application.CutCopyMode = 0
}

# ---- Overview sheet: b.md row (row 3) ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = $statusReady
$ovw.Range("F3").Value = $statusReady
$ovw.Range("G3").Value = "2016-08-18 16:36:15"

# ---- zh-cn sheet: b.md row (row 3) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $statusReady
Set-TextValue $zh.Range("F3") "False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-18 16:36:10"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: b.md row (row 3) ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $statusReady
Set-TextValue $de.Range("F3") "False"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-08-18 16:36:15"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 39.17
